# Auto-generated: apply updated profit figures per commit "chore: update Sheets via scheduled runner"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1978.9584
$ws.Range("I40").Value = 1266.6666
$ws.Range("J40").Value = 2406.3333
$ws.Range("K40").Value = 1266.6666
$ws.Range("L40").Value = 2406.3333
$ws.Range("M40").Value = -1091.6666
$ws.Range("N40").Value = -2756.3333

$ws.Range("H76").Value = 2965.1304
$ws.Range("I76").Value = 2751.3125
$ws.Range("J76").Value = 3453.8572
$ws.Range("K76").Value = 2751.3125
$ws.Range("L76").Value = 3453.8572
$ws.Range("M76").Value = -2436.3125
$ws.Range("N76").Value = -4083.8572

$ws.Range("H79").Value = 2965.1304
$ws.Range("I79").Value = 2751.3125
$ws.Range("J79").Value = 3453.8572
$ws.Range("K79").Value = 2751.3125
$ws.Range("L79").Value = 3453.8572
$ws.Range("M79").Value = -1659.3125
$ws.Range("N79").Value = -5637.8572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 680.6786
$ws.Range("I2").Value = 630
$ws.Range("J2").Value = 866.5
$ws.Range("K2").Value = 630
$ws.Range("L2").Value = 866.5
$ws.Range("M2").Value = -517
$ws.Range("N2").Value = -1092.5

$ws.Range("H32").Value = 3261.24
$ws.Range("I32").Value = 3052.8247
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 3052.8247
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -2765.8247
$ws.Range("N32").Value = -10574

$ws.Range("H63").Value = 3422.5217
$ws.Range("I63").Value = 2337.5405
$ws.Range("J63").Value = 7883
$ws.Range("K63").Value = 2337.5405
$ws.Range("L63").Value = 7883
$ws.Range("M63").Value = -1651.5405
$ws.Range("N63").Value = -9255

$ws.Range("H66").Value = 3422.5217
$ws.Range("I66").Value = 2337.5405
$ws.Range("J66").Value = 7883
$ws.Range("K66").Value = 11687.7025
$ws.Range("L66").Value = 39415
$ws.Range("M66").Value = -8255.702499999999
$ws.Range("N66").Value = -46279

$ws.Range("H88").Value = 1741.4445
$ws.Range("I88").Value = 1525
$ws.Range("J88").Value = 1914.6
$ws.Range("K88").Value = 1525
$ws.Range("L88").Value = 1914.6
$ws.Range("M88").Value = -1119
$ws.Range("N88").Value = -2726.6

$ws.Range("H91").Value = 1741.4445
$ws.Range("I91").Value = 1525
$ws.Range("J91").Value = 1914.6
$ws.Range("K91").Value = 1525
$ws.Range("L91").Value = 1914.6
$ws.Range("M91").Value = -121
$ws.Range("N91").Value = -4722.6

$ws.Range("H116").Value = 680.6786
$ws.Range("I116").Value = 630
$ws.Range("J116").Value = 866.5
$ws.Range("K116").Value = 630
$ws.Range("L116").Value = 866.5
$ws.Range("M116").Value = 1664
$ws.Range("N116").Value = -5454.5

$ws.Range("H122").Value = 1066.2
$ws.Range("I122").Value = 915
$ws.Range("J122").Value = 1387.5
$ws.Range("K122").Value = 2745
$ws.Range("L122").Value = 4162.5
$ws.Range("M122").Value = -295
$ws.Range("N122").Value = -9062.5

$ws.Range("H132").Value = 647170.9
$ws.Range("I132").Value = 1109804.8
$ws.Range("K132").Value = 3329414.4
$ws.Range("M132").Value = -3326884.4

$ws.Range("H133").Value = 42396.223
$ws.Range("J133").Value = 42396.223
$ws.Range("L133").Value = 42396.223
$ws.Range("N133").Value = -47456.223

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 680.6786
$ws.Range("I3").Value = 630
$ws.Range("J3").Value = 866.5
$ws.Range("K3").Value = 630
$ws.Range("L3").Value = 866.5
$ws.Range("M3").Value = -516
$ws.Range("N3").Value = -1094.5

$ws.Range("H105").Value = 2113.1428
$ws.Range("I105").Value = 1784.6154
$ws.Range("J105").Value = 2647
$ws.Range("K105").Value = 1784.6154
$ws.Range("L105").Value = 2647
$ws.Range("M105").Value = -37.61539999999991
$ws.Range("N105").Value = -6141

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 19949.166
$ws.Range("I62").Value = 3231.6667
$ws.Range("J62").Value = 36666.668
$ws.Range("K62").Value = 3231.6667
$ws.Range("L62").Value = 36666.668
$ws.Range("M62").Value = -2607.6667
$ws.Range("N62").Value = -37914.668

$ws.Range("H65").Value = 19949.166
$ws.Range("I65").Value = 3231.6667
$ws.Range("J65").Value = 36666.668
$ws.Range("K65").Value = 16158.3335
$ws.Range("L65").Value = 183333.34
$ws.Range("M65").Value = -13038.3335
$ws.Range("N65").Value = -189573.34

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2399.875
$ws.Range("I80").Value = 1000
$ws.Range("J80").Value = 2599.8572
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 7799.571599999999
$ws.Range("M80").Value = -2064
$ws.Range("N80").Value = -9671.571599999999

$ws.Range("H83").Value = 2399.875
$ws.Range("I83").Value = 1000
$ws.Range("J83").Value = 2599.8572
$ws.Range("K83").Value = 9000
$ws.Range("L83").Value = 23398.7148
$ws.Range("M83").Value = -4320
$ws.Range("N83").Value = -32758.7148

$ws.Range("H130").Value = 4298.2
$ws.Range("I130").Value = 5430
$ws.Range("J130").Value = 3543.6667
$ws.Range("K130").Value = 16290
$ws.Range("L130").Value = 10631.0001
$ws.Range("M130").Value = -11270
$ws.Range("N130").Value = -20671.0001

$ws.Range("H131").Value = 1029.9814
$ws.Range("I131").Value = 476.36365
$ws.Range("J131").Value = 1171.6046
$ws.Range("K131").Value = 1429.09095
$ws.Range("L131").Value = 3514.8138
$ws.Range("M131").Value = 3610.90905
$ws.Range("N131").Value = -13594.8138

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7145.8
$ws.Range("I70").Value = 7800.5386
$ws.Range("J70").Value = 6436.5
$ws.Range("K70").Value = 7800.5386
$ws.Range("L70").Value = 6436.5
$ws.Range("M70").Value = -7530.5386
$ws.Range("N70").Value = -6976.5

$ws.Range("H73").Value = 7145.8
$ws.Range("I73").Value = 7800.5386
$ws.Range("J73").Value = 6436.5
$ws.Range("K73").Value = 7800.5386
$ws.Range("L73").Value = 6436.5
$ws.Range("M73").Value = -6864.5386
$ws.Range("N73").Value = -8308.5

$ws.Range("H80").Value = 2783.7273
$ws.Range("I80").Value = 2670.1667
$ws.Range("J80").Value = 2920
$ws.Range("K80").Value = 2670.1667
$ws.Range("L80").Value = 2920
$ws.Range("M80").Value = -1672.1667
$ws.Range("N80").Value = -4916

$ws.Range("H83").Value = 2783.7273
$ws.Range("I83").Value = 2670.1667
$ws.Range("J83").Value = 2920
$ws.Range("K83").Value = 13350.8335
$ws.Range("L83").Value = 14600
$ws.Range("M83").Value = -8358.833500000001
$ws.Range("N83").Value = -24584

$ws.Range("H116").Value = 48371
$ws.Range("J116").Value = 48371
$ws.Range("L116").Value = 48371
$ws.Range("N116").Value = -57549

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 90911864
$ws.Range("I40").Value = 111113060
$ws.Range("K40").Value = 111113060
$ws.Range("M40").Value = -111112924

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 6292.25
$ws.Range("J104").Value = 6292.25
$ws.Range("L104").Value = 6292.25
$ws.Range("N104").Value = -13280.25

$ws.Range("H126").Value = 3566.875
$ws.Range("I126").Value = 3495.476
$ws.Range("J126").Value = 4066.6667
$ws.Range("K126").Value = 10486.428
$ws.Range("L126").Value = 12200.0001
$ws.Range("M126").Value = -8016.428
$ws.Range("N126").Value = -17140.0001
